$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45
$ws.Range("A45").Value = "2026/1/30"
$ws.Range("B45").Value = "2026/3/27"
$ws.Range("C45").Value = "第95期 秘寶 開放區域 極樂島 祕寶效果: 收穫盧恩文字10%(36%)翻倍"

# Row 46
$ws.Range("A46").Value = "2026/2/6"
$ws.Range("C46").Value = "第96期 混合紙飛機"
$ws.Range("B46").Value = "2026/4/3"

# Row 47
$ws.Range("A47").Value = "2026/2/13"
$ws.Range("B47").Value = "2026/4/10"
$ws.Range("C47").Value = "第97期 十一轉技能(星座)"

# Reflect the final selection/view state from the authored edit
[void]$ws.Range("C44").Select()
